$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 16.609
$ws.Range("A8").Value = -22.188
$ws.Range("A10").Value = -21.846
$ws.Range("A12").Value = -21.696
$ws.Range("E15").Value = 16.361
$ws.Range("A18").Value = -22.095
$ws.Range("E18").Value = 16.386
$ws.Range("E20").Value = 16.37
$ws.Range("E29").Value = 17.05
$ws.Range("E30").Value = 16.252
$ws.Range("E31").Value = 16.225
$ws.Range("A37").Value = -20.217
$ws.Range("E40").Value = 16.627
$ws.Range("E50").Value = 16.249
$ws.Range("A55").Value = -22.166
$ws.Range("A68").Value = -21.604
$ws.Range("E68").Value = 17.333
$ws.Range("E76").Value = 16.558
$ws.Range("A77").Value = -20.505
$ws.Range("A78").Value = -19.864
$ws.Range("A81").Value = -21.82
$ws.Range("A82").Value = -22.152
$ws.Range("E87").Value = 16.38
$ws.Range("E88").Value = 16.208
$ws.Range("E96").Value = 16.325
$ws.Range("E98").Value = 16.293
$ws.Range("E101").Value = 16.625
$ws.Range("E102").Value = 16.64
